$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new trailing columns ("nextkin" / "kinphone") mirroring the
# existing "payment" column (I) layout: bold/filled/bordered header in
# row 1, bordered (otherwise default) cells for the data rows below.

# Copy the formatting of column I (header style + data style) onto the
# new J and K columns first, so the freshly written header/data cells
# inherit the same cell styles used by the rest of the table.
$ws.Range("I1:I14").Copy()
$ws.Range("J1:J14").PasteSpecial(-4122)
$ws.Range("I1:I14").Copy()
$ws.Range("K1:K14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header text for the two new columns.
$ws.Range("J1").Value = "nextkin"
$ws.Range("K1").Value = "kinphone"

# Extend the duplicate-values conditional formatting that covered H2:H14
# so it also spans the two new columns (H2:J14).
$fcs = $ws.Range("H2:H14").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("H2:J14"))

# Match the author's final selection.
$null = $ws.Range("J8").Select()
